$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table grows from a 10x8 (mass 250..2500) x (coupling 0.25..3.5) grid to a
# 7x9 grid (mass 1000..2500) that also adds a coupling=1.8 row, i.e. every value
# effectively shifts left three columns and a new row is spliced in before the
# old "2.0" row.

# Stash the pristine bold/centered/bordered label style (used by the header row
# and the coupling column) in an unused scratch cell before any text rewrites -
# typing a numeric-looking quoted string later on flips a stray quote-prefix flag
# on the style actually used, so we keep a clean donor to paste back at the end.
$ws.Range("A2").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 1) Insert a new row at 6 so the coupling=1.8 row has room (rows 6-9 shift to 7-10)
$ws.Range("A6:K6").Insert()

# 2) The grid narrows from 10 mass points to 7, so columns I:J:K are no longer used
$ws.Range("I1:K10").Delete()

# 3) Header row (mediator masses, GeV) - stored as text, so force-quote the
#    numeric-looking literals to keep them text instead of numbers
$ws.Range("B1").Value = "'1000"
$ws.Range("C1").Value = "'1250"
$ws.Range("D1").Value = "'1500"
$ws.Range("E1").Value = "'1750"
$ws.Range("F1").Value = "'2000"
$ws.Range("G1").Value = "'2250"
$ws.Range("H1").Value = "'2500"

# 4) Coupling column (A2:A10) - also text
$ws.Range("A2").Value = "'0.25"
$ws.Range("A3").Value = "'0.5"
$ws.Range("A4").Value = "'1.0"
$ws.Range("A5").Value = "'1.5"
$ws.Range("A6").Value = "'1.8"
$ws.Range("A7").Value = "'2.0"
$ws.Range("A8").Value = "'2.5"
$ws.Range("A9").Value = "'3.0"
$ws.Range("A10").Value = "'3.5"

# 5) Cross-section matrix values (numeric)
$ws.Range("B2").Value = [double]"0.1974"
$ws.Range("C2").Value = [double]"0.0333"
$ws.Range("D2").Value = [double]"0.006502"
$ws.Range("E2").Value = [double]"0.001415"
$ws.Range("F2").Value = [double]"0.0003296"
$ws.Range("G2").Value = [double]"7.954e-05"
$ws.Range("H2").Value = [double]"1.98e-05"
$ws.Range("B3").Value = [double]"0.196"
$ws.Range("C3").Value = [double]"0.03336"
$ws.Range("D3").Value = [double]"0.006538"
$ws.Range("E3").Value = [double]"0.001414"
$ws.Range("F3").Value = [double]"0.0003285"
$ws.Range("G3").Value = [double]"7.96e-05"
$ws.Range("H3").Value = [double]"1.986e-05"
$ws.Range("B4").Value = [double]"0.1971"
$ws.Range("C4").Value = [double]"0.03331"
$ws.Range("D4").Value = [double]"0.006535"
$ws.Range("E4").Value = [double]"0.001416"
$ws.Range("F4").Value = [double]"0.0003302"
$ws.Range("G4").Value = [double]"7.931e-05"
$ws.Range("H4").Value = [double]"1.981e-05"
$ws.Range("B5").Value = [double]"0.1969"
$ws.Range("C5").Value = [double]"0.03326"
$ws.Range("D5").Value = [double]"0.006532"
$ws.Range("E5").Value = [double]"0.001418"
$ws.Range("F5").Value = [double]"0.0003282"
$ws.Range("G5").Value = [double]"7.983999999999999e-05"
$ws.Range("H5").Value = [double]"1.987e-05"
$ws.Range("B6").Value = [double]"0.1982092"
$ws.Range("C6").Value = [double]"0.03347837"
$ws.Range("D6").Value = [double]"0.006548228"
$ws.Range("E6").Value = [double]"0.0014225013"
$ws.Range("F6").Value = [double]"0.00033068328"
$ws.Range("G6").Value = [double]"7.9878115e-05"
$ws.Range("H6").Value = [double]"1.9879499e-05"
$ws.Range("B7").Value = [double]"0.1988"
$ws.Range("C7").Value = [double]"0.03358"
$ws.Range("D7").Value = [double]"0.006568"
$ws.Range("E7").Value = [double]"0.001436"
$ws.Range("F7").Value = [double]"0.0003313"
$ws.Range("G7").Value = [double]"7.972e-05"
$ws.Range("H7").Value = [double]"1.999e-05"
$ws.Range("B8").Value = [double]"0.2026"
$ws.Range("C8").Value = [double]"0.03419"
$ws.Range("D8").Value = [double]"0.006678"
$ws.Range("E8").Value = [double]"0.00145"
$ws.Range("F8").Value = [double]"0.0003353"
$ws.Range("G8").Value = [double]"8.069e-05"
$ws.Range("H8").Value = [double]"2.009e-05"
$ws.Range("B9").Value = [double]"0.2086"
$ws.Range("C9").Value = [double]"0.03522"
$ws.Range("D9").Value = [double]"0.006843"
$ws.Range("E9").Value = [double]"0.001482"
$ws.Range("F9").Value = [double]"0.0003419"
$ws.Range("G9").Value = [double]"8.231e-05"
$ws.Range("H9").Value = [double]"2.033e-05"
$ws.Range("B10").Value = [double]"0.219"
$ws.Range("C10").Value = [double]"0.03666"
$ws.Range("D10").Value = [double]"0.00716"
$ws.Range("E10").Value = [double]"0.001549"
$ws.Range("F10").Value = [double]"0.0003576"
$ws.Range("G10").Value = [double]"8.500000000000001e-05"
$ws.Range("H10").Value = [double]"2.092e-05"

# 6) Restore the bold/centered/bordered label style (clears the stray quote-prefix
#    flag picked up in steps 3-4) on every text cell, then drop the scratch cell
$ws.Range("Z1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)
$ws.Range("A2:A10").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("Z1").Clear()
